$d = $word.ActiveDocument

# --- Clean up the four log entries that previously had spell-check
# proofErr bookmarks splitting a single sentence across several runs
# (GitHub / wireframes / github). Running the identical visible text
# back through Find & Replace collapses each paragraph into a single
# plain run and removes the proofErr markers, matching the target.

$cleanups = @(
    "(19-04-21) Idag har laddat upp kodmappen till GitHub, skapat ASP kod som en grund med layout sidan, blivit klar med formuläret och skrivit vidare på målgruppsanalysen.",
    "(24-04-21) Jag är nu klar med min målgruppsanalys och har påbörjat wireframes.",
    "(27-04-21) Idag har jag jobbat vidare med wireframes och programmering.",
    "(28-04-21) Idag har jag kodat vidare och uppdaterat min github."
)

foreach ($text in $cleanups) {
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

# --- Add the new (04-05-21) log entry right after the (02-05-21) entry
# and before the trailing blank paragraphs at the end of the log.

# Locate the (02-05-21) paragraph through the Paragraphs collection so we
# get a live Paragraph object to anchor the insertion on.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("(02-05-21)")) {
        $target = $p
    }
}

$target.Range.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Range.InsertAfter("(04-05-21) Idag har jag fixat ThemeChanger javascriptet, ")

# Build the second run of the new paragraph as a temporary paragraph of
# its own (so the engine keeps it as a distinct <w:r>), then delete the
# paragraph mark between them to fold it back into one paragraph that
# holds two runs, exactly like the target markup.
$splitPoint = $d.Range($newPara.Range.End - 1, $newPara.Range.End - 1)
$splitPoint.InsertParagraphAfter()

$secondPara = $newPara.Next()
$secondPara.Range.InsertAfter("css och html.")

$paragraphMark = $d.Range($newPara.Range.End - 1, $newPara.Range.End)
$paragraphMark.Delete()
